$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full refreshed dataset for rows 2..11 (columns A..R).
# Columns A, B, C, E, F, G, I, O, R stay constant across all rows; D, H, J, K, L, M, N, P, Q vary.
$data = @(
    @(12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44326, 13, 100112021, "Ají", "Americana (o)", "Primera", 15, 30000, 30000, 30000, "`$/caja 25 kilos", "Provincia de Limarí", 1200, 25, "Hortaliza"),
    @(12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44193, 13, 100112021, "Ají", "Americana (o)", "Primera", 15, 46000, 46000, 46000, "`$/caja 15 kilos", "Provincia de Limarí", 3067, 15, "Hortaliza"),
    @(12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44425, 13, 100112021, "Ají", "Americana (o)", "Primera", 15, 75000, 75000, 75000, "`$/caja 25 kilos", "Provincia de Limarí", 3000, 25, "Hortaliza"),
    @(12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44343, 13, 100112021, "Ají", "Americana (o)", "Primera", 20, 36000, 36000, 36000, "`$/caja 25 kilos", "Provincia de Limarí", 1440, 25, "Hortaliza"),
    @(12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44221, 13, 100112021, "Ají", "Americana (o)", "Primera", 22, 24000, 25000, 24545, "`$/caja 25 kilos", "Provincia de Limarí", 982, 25, "Hortaliza"),
    @(12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44340, 13, 100112021, "Ají", "Americana (o)", "Primera", 15, 35000, 35000, 35000, "`$/caja 25 kilos", "Provincia de Limarí", 1400, 25, "Hortaliza"),
    @(12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44421, 13, 100112021, "Ají", "Americana (o)", "Primera", 15, 75000, 75000, 75000, "`$/caja 25 kilos", "Provincia de Limarí", 3000, 25, "Hortaliza"),
    @(12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44446, 13, 100112021, "Ají", "Americana (o)", "Primera", 5, 78000, 78000, 78000, "`$/caja 25 kilos", "Provincia de Limarí", 3120, 25, "Hortaliza"),
    @(12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44446, 13, 100112021, "Ají", "Inferno", "Primera", 4, 80000, 80000, 80000, "`$/caja 15 kilos", "Provincia de Limarí", 5333, 15, "Hortaliza"),
    @(12, "Mapocho Venta Directa de Santiago", "Metropolitana", 44319, 13, 100112021, "Ají", "Americana (o)", "Primera", 20, 30000, 30000, 30000, "`$/caja 25 kilos", "Provincia de Limarí", 1200, 25, "Hortaliza")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}

# Row 11 is a brand-new row; make sure its date cell (column D) carries the
# same date/time number format used by the rest of the column.
$ws.Range("D11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
